$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 16 (cohort_year=2019, period_index=6): num_customers 10 -> 11
$ws.Range("C16").Value = 11
$ws.Range("E16").Value = 0.005220692928334125

# Row 27 (cohort_year=2020, period_index=4): num_customers 52 -> 53
$ws.Range("C27").Value = 53
$ws.Range("E27").Value = 0.02353463587921847

# Row 34 (cohort_year=2021, period_index=2): num_customers 86 -> 87
$ws.Range("C34").Value = 87
$ws.Range("E34").Value = 0.03856382978723404

# Row 36 (cohort_year=2022, period_index=1): num_customers 142 -> 143
$ws.Range("C36").Value = 143
$ws.Range("E36").Value = 0.07409326424870466

# Row 37 (cohort_year=2023, period_index=0): num_customers/cohort_size 928 -> 937
$ws.Range("C37").Value = 937
$ws.Range("D37").Value = 937
$ws.Range("E37").Value = 1
